# Generate Report for Handoff
# - Updates the localization status text from "Handed back: in sync with en-US"
#   to "Ready for handoff" everywhere it appears (Overview!E2/F2, zh-cn!C2, de-de!C2).
# - Bumps the handoff-generation timestamps to the new run time.
# - Narrows the (now shorter) status column on every sheet to match.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$ws_overview.Range("E2").Value = "Ready for handoff"
$ws_overview.Range("F2").Value = "Ready for handoff"
$ws_zhcn.Range("C2").Value = "Ready for handoff"
$ws_dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
$ws_overview.Range("G2").Value = "2016-08-16 08:59:22"
$ws_zhcn.Range("H2").Value = "2016-08-16 08:59:16"
$ws_dede.Range("H2").Value = "2016-08-16 08:59:22"

# --- Column widths: status columns shrink from ~30 chars to ~17 chars ---
$ws_overview.Columns.Item(5).ColumnWidth = 16.33
$ws_overview.Columns.Item(6).ColumnWidth = 16.33
$ws_zhcn.Columns.Item(3).ColumnWidth = 16.33
$ws_dede.Columns.Item(3).ColumnWidth = 16.33
